$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fields")

# ---------------------------------------------------------------------------
# The "metric_ids" field (row 42 of the line_items.tsv block) is being split
# into three separate fields: jee3_metric_id, spar2_metric_id, hepr_metric_id.
# We duplicate row 42 twice (so it becomes three rows total: 42, 43, 44),
# pushing everything below down by two rows, then rewrite the content of the
# three rows with the new field definitions.
# ---------------------------------------------------------------------------

$ws.Rows.Item(42).Copy() | Out-Null
$ws.Rows.Item(42).Insert() | Out-Null
$ws.Rows.Item(42).Copy() | Out-Null
$ws.Rows.Item(42).Insert() | Out-Null

# Re-apply the original row's number formatting/borders to the two newly
# inserted rows (Insert after Copy brings the values but not the row-level
# formatting), then fix the row height which Excel also doesn't carry over
# automatically.
$ws.Range("A44:D44").Copy() | Out-Null
$ws.Range("A42:D42").PasteSpecial(-4122) | Out-Null
$ws.Range("A44:D44").Copy() | Out-Null
$ws.Range("A43:D43").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Rows.Item(42).RowHeight = 90
$ws.Rows.Item(43).RowHeight = 90
$ws.Rows.Item(44).RowHeight = 90

# Field names/definitions first (row 42 = jee3, 43 = spar2, 44 = hepr) ...
$ws.Range("B42").Value2 = "jee3_metric_id"
$ws.Range("C42").Value2 = "A unique ID associated with the specified metric of JEE 3.0"

$ws.Range("B43").Value2 = "spar2_metric_id"
$ws.Range("C43").Value2 = "A unique ID associated with the specified metric of SPAR 2.0"

$ws.Range("B44").Value2 = "hepr_metric_id"
$ws.Range("C44").Value2 = "A unique ID associated with the specified metric of Global Architecture for Health Emergency Preparedness, Response and Resilience (HEPR)"

# ... then the notes column, filled in HEPR / JEE / SPAR order.
$ws.Range("D44").Value2 = "Where a single line-item corresponding to more than one metric from a given framework (e.g. JEE or SPAR), efforts where made to disambiguate to the extent possible in order to assign each line item to a maximum of one specific metric per framework. For example, a given cost will not be mapped to more than one metric of the JEE. Instead, each line item was mapped to the single most relevant metric per framework. For more information on the HEPR framework, please see WHO documentation https://www.who.int/emergencies/operations/universal-health---preparedness-review"

$ws.Range("D42").Value2 = "Where a single line-item corresponding to more than one metric from a given framework (e.g. JEE or SPAR), efforts where made to disambiguate to the extent possible in order to assign each line item to a maximum of one specific metric per framework. For example, a given cost will not be mapped to more than one metric of the JEE. Instead, each line item was mapped to the single most relevant metric per framework. For more information on the JEE, please see WHO documentation https://www.who.int/publications/i/item/9789240051980"

$ws.Range("D43").Value2 = "Where a single line-item corresponding to more than one metric from a given framework (e.g. JEE or SPAR), efforts where made to disambiguate to the extent possible in order to assign each line item to a maximum of one specific metric per framework. For example, a given cost will not be mapped to more than one metric of the JEE. Instead, each line item was mapped to the single most relevant metric per framework. For more information on the SPAR, please see WHO documentation https://www.who.int/emergencies/operations/international-health-regulations-monitoring-evaluation-framework/states-parties-self-assessment-annual-reporting"

Write-Host "Done"
